# Update "Related", "Visited" and "Final Selected" counts for strategy E1 (row 2),
# and recompute Precision / Recall / F-Measure / Final Precision / Final Recall /
# Final F-Measure for rows 2-8, reflecting the corrected matriz de citacao numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (E1) - counts changed
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 43
$ws.Range("E2").Value = 4

# Row 2 (E1) - recomputed ratios
$ws.Range("F2").Value = 0.1395348837209302
$ws.Range("G2").Value = 0.1176470588235294
$ws.Range("H2").Value = 0.1276595744680851
$ws.Range("I2").Value = 0.09302325581395349
$ws.Range("J2").Value = 0.1333333333333333
$ws.Range("K2").Value = 0.1095890410958904

# Row 3 (E2) - recomputed ratios
$ws.Range("F3").Value = 0.02722904431393487
$ws.Range("H3").Value = 0.05301455301455301
$ws.Range("I3").Value = 0.01601708489054992
$ws.Range("K3").Value = 0.03152916447714135

# Row 4 (E3) - recomputed ratios
$ws.Range("F4").Value = 0.0339943342776204
$ws.Range("G4").Value = 0.7058823529411765
$ws.Range("H4").Value = 0.06486486486486487
$ws.Range("I4").Value = 0.02077431539187913
$ws.Range("J4").Value = 0.7333333333333333
$ws.Range("K4").Value = 0.04040404040404041

# Row 5 (E4) - recomputed ratios
$ws.Range("F5").Value = 0.03747870528109029
$ws.Range("G5").Value = 0.8627450980392157
$ws.Range("H5").Value = 0.07183673469387755
$ws.Range("I5").Value = 0.02129471890971039
$ws.Range("J5").Value = 0.8333333333333334
$ws.Range("K5").Value = 0.04152823920265781

# Row 6 (E5) - recomputed ratios
$ws.Range("F6").Value = 0.06506849315068493
$ws.Range("G6").Value = 0.3725490196078431
$ws.Range("H6").Value = 0.1107871720116618
$ws.Range("I6").Value = 0.04794520547945205
$ws.Range("J6").Value = 0.4666666666666667
$ws.Range("K6").Value = 0.08695652173913043

# Row 7 (E6) - recomputed ratios
$ws.Range("F7").Value = 0.06194690265486726
$ws.Range("G7").Value = 0.6862745098039216
$ws.Range("H7").Value = 0.1136363636363636
$ws.Range("I7").Value = 0.03539823008849557
$ws.Range("J7").Value = 0.6666666666666666
$ws.Range("K7").Value = 0.06722689075630252

# Row 8 (E7) - recomputed ratios
$ws.Range("F8").Value = 0.05811138014527845
$ws.Range("G8").Value = 0.4705882352941176
$ws.Range("H8").Value = 0.103448275862069
$ws.Range("I8").Value = 0.03631961259079903
$ws.Range("K8").Value = 0.06772009029345372
